$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two more rows shaped like row 3 (rows 4 and 5), then one more
# shaped like row 2 (row 6) - same data, reusing the existing shared
# strings / styles via a direct range-to-range copy (keeps cell styles).
$ws.Range("A3:AH3").Copy($ws.Range("A4:AH4"))
$ws.Range("A3:AH3").Copy($ws.Range("A5:AH5"))
$ws.Range("A2:AH2").Copy($ws.Range("A6:AH6"))

# Move the selection / scroll position like the UI would after this work.
[void]$ws.Range("F12").Select()

Write-Output "done"
